# Update the Files upload section headers/data on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: meterId -> userId (B1/C1 remain unchanged)
$ws.Range("A1").Value = "userId"

# Data rows: replace meter IDs with simple sequential userIds
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Update the active selection to A2
$ws.Range("A2").Select()

$wb.Save()
